$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("Docente(s) Responsável(eis) ", $false, $false, $false, $false, $false,
                $true, 1, $false, "", 0)

$r.InsertParagraphAfter()

# $r now spans the found text (pre-insert); move to the newly created paragraph.
$newPar = $r.Paragraphs(1).Next()
$newPar.Range.Text = "1814052 - Silvio Silverio da Silva"
$newPar.Style = "ListBullet"
